{"js": "const body = context.document.body;\nconst paras = body.paragraphs;\nparas.load(\"items\");\nawait context.sync();\n\n// --- Paragraph 1 (index 0): update date, add line break + new title line ---\n// Built via insertOoxml so the <w:br/> + second <w:t> land inside the SAME\n// run as the first <w:t>, matching the canonical OOXML produced by Word.\nconst p0 = paras.items[0];\nconst p0Ooxml =\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body><w:p><w:pPr><w:pStyle w:val=\"Normal\"/></w:pPr><w:r>' +\n  '<w:t>\\u26A1\\uFE0F\\uD83D\\uDE80\u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d9\u05d5\u05de\u05d9 \u05e9\u05dc \u05de\u05d9\u05d9\u05e7 -08.11.24: \\u26A1\\uFE0F\\uD83D\\uDE80</w:t>' +\n  '<w:br/>' +\n  '<w:t>Occam\\u2019s Razor for Self Supervised Learning: What is Sufficient to Learn Good Representations?</w:t>' +\n  '</w:r></w:p></w:body></w:document>' +\n  '</pkg:xmlData></pkg:part></pkg:package>';\np0.insertOoxml(p0Ooxml, \"Replace\");\n\n// --- Paragraphs 2-6 (index 1-5): full text replacement ---\nconst p1 = paras.items[1];\np1.insertText(\n  \"\u05e1\u05e7\u05d9\u05e8\u05d4 \u05e7\u05e6\u05e8\u05d4 \u05e9\u05dc \u05de\u05d0\u05de\u05e8 \u05d4\u05de\u05e6\u05d9\u05e2 \u05d2\u05d9\u05e9\u05d4 \u05d7\u05d3\u05e9\u05d4 \u05dc\u05dc\u05de\u05d9\u05d3\u05d4 self-supervised \u05d0\u05d5 SSL \u05d1\u05e7\u05e6\u05e8\u05d4. \u05d0\u05d6\u05db\u05d9\u05e8 \u05db\u05d9 \u05e9\u05d9\u05d8\u05ea SSL \u05de\u05e0\u05d9\u05d7\u05d4 \u05e9\u05d9\u05e9 \u05dc\u05e0\u05d5 \u05d3\u05d0\u05d8\u05d4 \u05dc\u05d0 \u05de\u05ea\u05d5\u05d9\u05d2 \u05d5\u05de\u05ea\u05e8\u05d8\u05e0\u05d5 \u05dc\u05d0\u05de\u05df \u05de\u05d5\u05d3\u05dc \u05de\u05e1\u05d5\u05d2\u05dc \u05dc\u05d4\u05e4\u05d9\u05e7 \u05d9\u05d9\u05e6\u05d5\u05d2 \u05d7\u05d6\u05e7 \u05e9\u05dc \u05d3\u05d0\u05d8\u05d4. \u05de\u05d4 \u05d6\u05d4 \u05d9\u05d9\u05e6\u05d5\u05d2 \u05d7\u05d6\u05e7 \u05e9\u05dc \u05d3\u05d0\u05d8\u05d4, \u05d0\u05ea\u05dd \u05e9\u05d5\u05d0\u05dc\u05d9\u05dd? \u05d1\u05d3\u05f4\u05db \u05d4\u05db\u05d5\u05d5\u05e0\u05d4 \u05dc\u05db\u05d6\u05d4 \u05e9\u05e0\u05d9\u05ea\u05df \u05dc\u05de\u05e0\u05e3 \u05d0\u05d5\u05ea\u05d5 \u05d1\u05e6\u05d5\u05e8\u05d4 \u05e7\u05dc\u05d4 (\u05e0\u05d2\u05d9\u05d3 \u05e8\u05e7 \u05e2\u05dd \u05ea\u05d5\u05e1\u05e4\u05ea \u05e9\u05dc \u05e9\u05db\u05d1\u05d4 \u05dc\u05d9\u05e0\u05d0\u05e8\u05d9\u05ea) \u05dc\u05d1\u05e0\u05d9\u05d9\u05ea \u05de\u05e1\u05d5\u05d5\u05d2 \u05d1\u05e2\u05dc \u05d1\u05d9\u05e6\u05d5\u05e2\u05d9\u05dd \u05d8\u05d5\u05d1\u05d9\u05dd. \",\n  \"Replace\"\n);\n\nconst p2 = paras.items[2];\np2.insertText(\n  \"\u05db\u05dc\u05d5\u05de\u05e8 \u05db\u05d6\u05d4 \u05e9\u05d9\u05d5\u05d3\u05e2 \u05dc\u05d4\u05e4\u05e8\u05d9\u05d3 \u05d1\u05d9\u05df \u05d4\u05e7\u05d8\u05d2\u05d5\u05e8\u05d9\u05d5\u05ea \u05d4\u05e9\u05d5\u05e0\u05d5\u05ea \u05e9\u05dc \u05d3\u05d0\u05d8\u05d4 \u05d1\u05dc\u05d9 \u05dc\u05d3\u05e2\u05ea \u05d0\u05d5\u05ea\u05df \u05d1\u05e6\u05d5\u05e8\u05d4 \u05de\u05e4\u05d5\u05e8\u05e9\u05ea (\u05dc\u05de\u05e9\u05dc \u05d0\u05e0\u05d5 \u05d9\u05db\u05d5\u05dc\u05d9\u05dd \u05dc\u05d0\u05de\u05df \u05de\u05d5\u05d3\u05dc \u05d1\u05e6\u05d5\u05e8\u05ea SSL \u05e2\u05dc \u05d4\u05ea\u05de\u05d5\u05e0\u05d5\u05ea \u05e9\u05dc ImageNet \u05d1\u05dc\u05d9 \u05dc\u05d4\u05e9\u05ea\u05de\u05e9 \u05d1\u05ea\u05d9\u05d5\u05d2\u05d9\u05dd \u05d5\u05d0\u05d6 \u05dc\u05d1\u05d3\u05d5\u05e7 \u05d4\u05d0\u05dd \u05d4\u05de\u05d5\u05d3\u05dc \u05d4\u05e6\u05dc\u05d9\u05d7 \u05dc\u05dc\u05de\u05d5\u05d3 \u05dc\u05d4\u05e4\u05e8\u05d9\u05d3 \u05d1\u05d9\u05df \u05d4\u05e7\u05d8\u05d2\u05d5\u05e8\u05d9\u05d5\u05ea \u05d4\u05e9\u05d5\u05e0\u05d5\u05ea).\",\n  \"Replace\"\n);\n\n// NOTE: this paragraph's run originally had xml:space=\"preserve\" (its old\n// text ended with a trailing space). insertText(...,\"Replace\") keeps that\n// attribute stuck on the run even though the new text has no leading/\n// trailing whitespace, so insertOoxml is used instead for an exact rebuild\n// of the run (matching the canonical OOXML, which has no xml:space here).\nconst p3 = paras.items[3];\nconst p3Ooxml =\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body><w:p><w:pPr><w:pStyle w:val=\"Normal\"/></w:pPr><w:r>' +\n  '<w:t>\u05d1\u05d3\u05f4\u05db SSL \u05de\u05d1\u05d5\u05e6\u05e2 \u05e2\u05dd \u05e9\u05d9\u05d8\u05d5\u05ea \u05e9\u05dc \u05dc\u05de\u05d9\u05d3\u05d4 \u05e0\u05d9\u05d2\u05d5\u05d3\u05d9\u05ea (contrastive learning) \u05db\u05d0\u05e9\u05e8 \u05de\u05d8\u05e8\u05ea\u05d5 \u05de\u05d0\u05d5\u05d3 \u05d1\u05d2\u05d3\u05d5\u05dc \u05d4\u05d9\u05d0 \u05dc\u05e7\u05e8\u05d1 \u05d9\u05d9\u05e6\u05d5\u05d2\u05d9\u05dd \u05e9\u05dc \u05e4\u05d9\u05e1\u05d5\u05ea \u05d3\u05d0\u05d8\u05d4 \u05d3\u05d5\u05de\u05d5\u05ea (\u05d7\u05d9\u05d5\u05d1\u05d9\u05d5\u05ea) \u05d5\u05dc\u05d4\u05e8\u05d7\u05d9\u05e7 \u05d0\u05ea \u05d4\u05d9\u05d9\u05e6\u05d5\u05d2\u05d9\u05dd \u05e9\u05dc \u05e4\u05d9\u05e1\u05d5\u05ea \u05d3\u05d0\u05d8\u05d4 \u05dc\u05d0 \u05d3\u05d5\u05de\u05d5\u05ea (\u05e9\u05dc\u05d9\u05dc\u05d9\u05d5\u05ea). \u05dc\u05e8\u05d5\u05d1 \u05d6\u05d5\u05d2\u05d5\u05ea \u05d7\u05d9\u05d5\u05d1\u05d9\u05d9\u05dd \u05e0\u05d1\u05d7\u05e8\u05d9\u05dd \u05d1\u05ea\u05d5\u05e8 \u05d0\u05d5\u05d2\u05de\u05e0\u05d8\u05e6\u05d9\u05d5\u05ea \u05e9\u05d5\u05e0\u05d5\u05ea \u05e9\u05dc \u05d0\u05d5\u05ea\u05d4 \u05d4\u05d3\u05d5\u05d2\u05de\u05d0 \u05db\u05d0\u05e9\u05e8 \u05d4\u05d6\u05d5\u05d2\u05d5\u05ea \u05d4\u05e9\u05dc\u05d9\u05dc\u05d9\u05d9\u05dd \u05d4\u05df \u05d3\u05d5\u05d2\u05de\u05d0\u05d5\u05ea \u05e9\u05e0\u05d1\u05d7\u05e8\u05d5\u05ea \u05d1\u05d0\u05e7\u05e8\u05d0\u05d9. \u05e9\u05d9\u05d8\u05d5\u05ea \u05db\u05d0\u05dc\u05d5 \u05e0\u05d7\u05dc\u05d5 \u05d4\u05e6\u05dc\u05d7\u05d4 \u05d3\u05d9 \u05d2\u05d3\u05d5\u05dc\u05d4 \u05d0\u05d1\u05dc \u05d3\u05e8\u05e9\u05d5 \u05d3\u05d0\u05d8\u05d4\u05e1\u05d8\u05d9\u05dd \u05de\u05d0\u05d5\u05d3 \u05d2\u05d3\u05d5\u05dc\u05d9\u05dd \u05d5\u05d2\u05dd \u05de\u05e9\u05d0\u05d1\u05d9 \u05d0\u05d9\u05de\u05d5\u05df \u05d3\u05d9 \u05de\u05e9\u05de\u05e2\u05d5\u05ea\u05d9\u05d9\u05dd (\u05db\u05d9 \u05e0\u05d3\u05e8\u05e9 \u05e9\u05dd \u05d2\u05d5\u05d3\u05dc \u05d1\u05d0\u05e5\\' \u05d3\u05d9 \u05d2\u05d3\u05d5\u05dc \u05db\u05d3\u05d9 \u05e9\u05d4\u05e9\u05d9\u05d8\u05d4 \u05ea\u05e2\u05d1\u05d5\u05d3 \u05d8\u05d5\u05d1).</w:t>' +\n  '</w:r></w:p></w:body></w:document>' +\n  '</pkg:xmlData></pkg:part></pkg:package>';\np3.insertOoxml(p3Ooxml, \"Replace\");\n\nconst p4 = paras.items[4];\np4.insertText(\n  \"\u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05de\u05e1\u05d5\u05e7\u05e8 \u05de\u05e6\u05d9\u05e2 \u05e9\u05d9\u05d8\u05d4 \u05de\u05d0\u05d5\u05d3 \u05e4\u05e9\u05d5\u05d8\u05d4 \u05d5\u05d0\u05d9\u05e0\u05d8\u05d5\u05d0\u05d8\u05d9\u05d1\u05d9\u05ea \u05dc-SSL(\u05ea\u05e2\u05e8 \u05d0\u05d5\u05e7\u05dd). \u05d1\u05de\u05e7\u05d5\u05dd \u05dc\u05e2\u05d1\u05d5\u05d3 \u05e2\u05dd \u05d4\u05d9\u05d9\u05e6\u05d5\u05d2\u05d9\u05dd \u05d4\u05de\u05d0\u05de\u05e8 \u05de\u05d0\u05de\u05df \u05de\u05d5\u05d3\u05dc \u05dc\u05d7\u05d6\u05d5\u05ea \u05d0\u05ea \u05d4\u05de\u05e1\u05e4\u05e8 \u05e9\u05dc \u05d4\u05d3\u05d5\u05d2\u05de\u05d0 \u05d1\u05d3\u05d0\u05d8\u05d4\u05e1\u05d8. \u05db\u05dc\u05d5\u05de\u05e8 \u05d0\u05dd \u05d9\u05e9 \u05dc\u05e0\u05d5 1000 \u05d3\u05d5\u05d2\u05de\u05d0\u05d5\u05ea \u05de\u05d4\u05d3\u05d0\u05d8\u05d4\u05e1\u05d8 \u05d9\u05e9 \u05dc\u05e0\u05d5 1000 \u05e7\u05d8\u05d2\u05d5\u05e8\u05d9\u05d5\u05ea \u05d5\u05de\u05d8\u05e8\u05ea\u05e0\u05d5 \u05dc\u05d7\u05d6\u05d5\u05ea \u05e7\u05d8\u05d2\u05d5\u05e8\u05d9\u05d4 \u05e9\u05dc \u05d3\u05d5\u05d2\u05de\u05d0 \u05de\u05d4\u05d9\u05d9\u05e6\u05d5\u05d2 \u05d4\u05dc\u05d8\u05e0\u05d8\u05d9 \u05e9\u05dc\u05d4 (\u05d4\u05de\u05d5\u05e4\u05e7 \u05e2\u05dc \u05d9\u05d3\u05d9 \u05d4\u05de\u05d5\u05d3\u05dc \u05d4\u05de\u05d0\u05d5\u05de\u05df). \u05db\u05dc\u05d5\u05de\u05e8 \u05d0\u05d7\u05e8\u05d9 \u05d4\u05e9\u05db\u05d1\u05d4 \u05d4\u05d0\u05d7\u05e8\u05d5\u05e0\u05d4 \u05e9\u05dc \u05d4\u05de\u05d5\u05d3\u05dc \u05de\u05d5\u05e1\u05d9\u05e4\u05d9\u05dd \u05e9\u05db\u05d1\u05d4 \u05e2\u05dd \u05de\u05d8\u05e8\u05d9\u05e6\u05d4 \u05d4\u05de\u05de\u05e4\u05d4 \u05d0\u05ea \u05d4\u05d9\u05d9\u05e6\u05d5\u05d2 \u05dc\u05e7\u05d8\u05d2\u05d5\u05e8\u05d9\u05d5\u05ea (\u05db\u05dc\u05d5\u05de\u05e8 \u05d4\u05de\u05e1\u05e4\u05e8\u05d9\u05dd \u05d4\u05e1\u05d9\u05d3\u05d5\u05e8\u05d9\u05d9\u05dd \u05e9\u05dc \u05d4\u05d3\u05d5\u05d2\u05de\u05d0\u05d5\u05ea). \u05d5\u05d1\u05e1\u05d5\u05e3 \u05e9\u05dc \u05dc\u05d5\u05e1 cross-entropy \u05d4\u05e1\u05d8\u05e0\u05d3\u05e8\u05d8\u05d9.\",\n  \"Replace\"\n);\n\nconst p5 = paras.items[5];\np5.insertText(\n  \"\u05d0\u05d6 \u05d4\u05de\u05d0\u05de\u05e8 \u05de\u05d5\u05db\u05d9\u05d7 \u05e9\u05d4\u05e9\u05d9\u05d8\u05d4 \u05e2\u05d5\u05d1\u05d3\u05ea \u05dc\u05d0 \u05e8\u05e2 \u05dc\u05d3\u05d0\u05d8\u05d4\u05e1\u05d8\u05d9\u05dd \u05d9\u05d7\u05e1\u05d9\u05ea \u05dc\u05d0 \u05d2\u05d3\u05d5\u05dc\u05d9\u05dd (\u05de\u05e2\u05e0\u05d9\u05d9\u05df \u05d0\u05d9\u05da \u05d6\u05d4 \u05d9\u05e2\u05d1\u05d5\u05d3 \u05dc\u05d3\u05d0\u05d8\u05d4\u05e1\u05d8 \u05d1\u05d2\u05d5\u05d3\u05dc 10 \u05de\u05d9\u05dc\u05d9\u05d5\u05df). \u05db\u05de\u05d5\u05d1\u05df \u05d9\u05e9 \u05db\u05de\u05d4 \u05d8\u05e8\u05d9\u05e7\u05d9\u05dd \u05d1\u05d0\u05d9\u05de\u05d5\u05df \u05db\u05de\u05d5 soft labels \u05d0\u05d1\u05dc \u05d1\u05d2\u05d3\u05d5\u05dc \u05d4\u05e8\u05e2\u05d9\u05d5\u05df \u05d3\u05d9 \u05e0\u05d7\u05de\u05d3. \",\n  \"Replace\"\n);\n\n// --- Paragraphs 7-9 (index 6-8): entirely removed ---\nparas.items[6].delete();\nparas.items[7].delete();\nparas.items[8].delete();\n\n// --- Paragraph 10 (now last, index 9): replace URL text ---\nconst pUrl = paras.items[9];\npUrl.insertText(\"https://arxiv.org/pdf/2406.10743\", \"Replace\");\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Helper: replace a paragraph's text without the old run's xml:space\n# attribute \"sticking\" to the new text (InsertParagraphBefore() creates a\n# brand-new run to hold the replacement, then the original paragraph -\n# now shifted one slot later - is deleted).\nfunction Set-ParagraphText($index, $text) {\n    $old = $d.Paragraphs($index).Range\n    $old.InsertParagraphBefore()\n    $d.Paragraphs($index).Range.Text = $text\n    $d.Paragraphs($index + 1).Range.Delete()\n}\n\n# --- Paragraph 1: fix the date, then append a manual line break + new title ---\n$find = $d.Content.Find\n$find.Text = \"09.11.24\"\n$find.Replacement.Text = \"08.11.24\"\n$find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2)\n\n$p1 = $d.Paragraphs(1).Range\n$p1.InsertAfter([char]11 + \"Occam\" + [char]0x2019 + \"s Razor for Self Supervised Learning: What is Sufficient to Learn Good Representations?\")\n\n# --- Paragraphs 2-6: full text replacement ---\nSet-ParagraphText 2 \"\u05e1\u05e7\u05d9\u05e8\u05d4 \u05e7\u05e6\u05e8\u05d4 \u05e9\u05dc \u05de\u05d0\u05de\u05e8 \u05d4\u05de\u05e6\u05d9\u05e2 \u05d2\u05d9\u05e9\u05d4 \u05d7\u05d3\u05e9\u05d4 \u05dc\u05dc\u05de\u05d9\u05d3\u05d4 self-supervised \u05d0\u05d5 SSL \u05d1\u05e7\u05e6\u05e8\u05d4. \u05d0\u05d6\u05db\u05d9\u05e8 \u05db\u05d9 \u05e9\u05d9\u05d8\u05ea SSL \u05de\u05e0\u05d9\u05d7\u05d4 \u05e9\u05d9\u05e9 \u05dc\u05e0\u05d5 \u05d3\u05d0\u05d8\u05d4 \u05dc\u05d0 \u05de\u05ea\u05d5\u05d9\u05d2 \u05d5\u05de\u05ea\u05e8\u05d8\u05e0\u05d5 \u05dc\u05d0\u05de\u05df \u05de\u05d5\u05d3\u05dc \u05de\u05e1\u05d5\u05d2\u05dc \u05dc\u05d4\u05e4\u05d9\u05e7 \u05d9\u05d9\u05e6\u05d5\u05d2 \u05d7\u05d6\u05e7 \u05e9\u05dc \u05d3\u05d0\u05d8\u05d4. \u05de\u05d4 \u05d6\u05d4 \u05d9\u05d9\u05e6\u05d5\u05d2 \u05d7\u05d6\u05e7 \u05e9\u05dc \u05d3\u05d0\u05d8\u05d4, \u05d0\u05ea\u05dd \u05e9\u05d5\u05d0\u05dc\u05d9\u05dd? \u05d1\u05d3\u05f4\u05db \u05d4\u05db\u05d5\u05d5\u05e0\u05d4 \u05dc\u05db\u05d6\u05d4 \u05e9\u05e0\u05d9\u05ea\u05df \u05dc\u05de\u05e0\u05e3 \u05d0\u05d5\u05ea\u05d5 \u05d1\u05e6\u05d5\u05e8\u05d4 \u05e7\u05dc\u05d4 (\u05e0\u05d2\u05d9\u05d3 \u05e8\u05e7 \u05e2\u05dd \u05ea\u05d5\u05e1\u05e4\u05ea \u05e9\u05dc \u05e9\u05db\u05d1\u05d4 \u05dc\u05d9\u05e0\u05d0\u05e8\u05d9\u05ea) \u05dc\u05d1\u05e0\u05d9\u05d9\u05ea \u05de\u05e1\u05d5\u05d5\u05d2 \u05d1\u05e2\u05dc \u05d1\u05d9\u05e6\u05d5\u05e2\u05d9\u05dd \u05d8\u05d5\u05d1\u05d9\u05dd. \"\n\nSet-ParagraphText 3 \"\u05db\u05dc\u05d5\u05de\u05e8 \u05db\u05d6\u05d4 \u05e9\u05d9\u05d5\u05d3\u05e2 \u05dc\u05d4\u05e4\u05e8\u05d9\u05d3 \u05d1\u05d9\u05df \u05d4\u05e7\u05d8\u05d2\u05d5\u05e8\u05d9\u05d5\u05ea \u05d4\u05e9\u05d5\u05e0\u05d5\u05ea \u05e9\u05dc \u05d3\u05d0\u05d8\u05d4 \u05d1\u05dc\u05d9 \u05dc\u05d3\u05e2\u05ea \u05d0\u05d5\u05ea\u05df \u05d1\u05e6\u05d5\u05e8\u05d4 \u05de\u05e4\u05d5\u05e8\u05e9\u05ea (\u05dc\u05de\u05e9\u05dc \u05d0\u05e0\u05d5 \u05d9\u05db\u05d5\u05dc\u05d9\u05dd \u05dc\u05d0\u05de\u05df \u05de\u05d5\u05d3\u05dc \u05d1\u05e6\u05d5\u05e8\u05ea SSL \u05e2\u05dc \u05d4\u05ea\u05de\u05d5\u05e0\u05d5\u05ea \u05e9\u05dc ImageNet \u05d1\u05dc\u05d9 \u05dc\u05d4\u05e9\u05ea\u05de\u05e9 \u05d1\u05ea\u05d9\u05d5\u05d2\u05d9\u05dd \u05d5\u05d0\u05d6 \u05dc\u05d1\u05d3\u05d5\u05e7 \u05d4\u05d0\u05dd \u05d4\u05de\u05d5\u05d3\u05dc \u05d4\u05e6\u05dc\u05d9\u05d7 \u05dc\u05dc\u05de\u05d5\u05d3 \u05dc\u05d4\u05e4\u05e8\u05d9\u05d3 \u05d1\u05d9\u05df \u05d4\u05e7\u05d8\u05d2\u05d5\u05e8\u05d9\u05d5\u05ea \u05d4\u05e9\u05d5\u05e0\u05d5\u05ea).\"\n\nSet-ParagraphText 4 \"\u05d1\u05d3\u05f4\u05db SSL \u05de\u05d1\u05d5\u05e6\u05e2 \u05e2\u05dd \u05e9\u05d9\u05d8\u05d5\u05ea \u05e9\u05dc \u05dc\u05de\u05d9\u05d3\u05d4 \u05e0\u05d9\u05d2\u05d5\u05d3\u05d9\u05ea (contrastive learning) \u05db\u05d0\u05e9\u05e8 \u05de\u05d8\u05e8\u05ea\u05d5 \u05de\u05d0\u05d5\u05d3 \u05d1\u05d2\u05d3\u05d5\u05dc \u05d4\u05d9\u05d0 \u05dc\u05e7\u05e8\u05d1 \u05d9\u05d9\u05e6\u05d5\u05d2\u05d9\u05dd \u05e9\u05dc \u05e4\u05d9\u05e1\u05d5\u05ea \u05d3\u05d0\u05d8\u05d4 \u05d3\u05d5\u05de\u05d5\u05ea (\u05d7\u05d9\u05d5\u05d1\u05d9\u05d5\u05ea) \u05d5\u05dc\u05d4\u05e8\u05d7\u05d9\u05e7 \u05d0\u05ea \u05d4\u05d9\u05d9\u05e6\u05d5\u05d2\u05d9\u05dd \u05e9\u05dc \u05e4\u05d9\u05e1\u05d5\u05ea \u05d3\u05d0\u05d8\u05d4 \u05dc\u05d0 \u05d3\u05d5\u05de\u05d5\u05ea (\u05e9\u05dc\u05d9\u05dc\u05d9\u05d5\u05ea). \u05dc\u05e8\u05d5\u05d1 \u05d6\u05d5\u05d2\u05d5\u05ea \u05d7\u05d9\u05d5\u05d1\u05d9\u05d9\u05dd \u05e0\u05d1\u05d7\u05e8\u05d9\u05dd \u05d1\u05ea\u05d5\u05e8 \u05d0\u05d5\u05d2\u05de\u05e0\u05d8\u05e6\u05d9\u05d5\u05ea \u05e9\u05d5\u05e0\u05d5\u05ea \u05e9\u05dc \u05d0\u05d5\u05ea\u05d4 \u05d4\u05d3\u05d5\u05d2\u05de\u05d0 \u05db\u05d0\u05e9\u05e8 \u05d4\u05d6\u05d5\u05d2\u05d5\u05ea \u05d4\u05e9\u05dc\u05d9\u05dc\u05d9\u05d9\u05dd \u05d4\u05df \u05d3\u05d5\u05d2\u05de\u05d0\u05d5\u05ea \u05e9\u05e0\u05d1\u05d7\u05e8\u05d5\u05ea \u05d1\u05d0\u05e7\u05e8\u05d0\u05d9. \u05e9\u05d9\u05d8\u05d5\u05ea \u05db\u05d0\u05dc\u05d5 \u05e0\u05d7\u05dc\u05d5 \u05d4\u05e6\u05dc\u05d7\u05d4 \u05d3\u05d9 \u05d2\u05d3\u05d5\u05dc\u05d4 \u05d0\u05d1\u05dc \u05d3\u05e8\u05e9\u05d5 \u05d3\u05d0\u05d8\u05d4\u05e1\u05d8\u05d9\u05dd \u05de\u05d0\u05d5\u05d3 \u05d2\u05d3\u05d5\u05dc\u05d9\u05dd \u05d5\u05d2\u05dd \u05de\u05e9\u05d0\u05d1\u05d9 \u05d0\u05d9\u05de\u05d5\u05df \u05d3\u05d9 \u05de\u05e9\u05de\u05e2\u05d5\u05ea\u05d9\u05d9\u05dd (\u05db\u05d9 \u05e0\u05d3\u05e8\u05e9 \u05e9\u05dd \u05d2\u05d5\u05d3\u05dc \u05d1\u05d0\u05e5' \u05d3\u05d9 \u05d2\u05d3\u05d5\u05dc \u05db\u05d3\u05d9 \u05e9\u05d4\u05e9\u05d9\u05d8\u05d4 \u05ea\u05e2\u05d1\u05d5\u05d3 \u05d8\u05d5\u05d1).\"\n\nSet-ParagraphText 5 \"\u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05de\u05e1\u05d5\u05e7\u05e8 \u05de\u05e6\u05d9\u05e2 \u05e9\u05d9\u05d8\u05d4 \u05de\u05d0\u05d5\u05d3 \u05e4\u05e9\u05d5\u05d8\u05d4 \u05d5\u05d0\u05d9\u05e0\u05d8\u05d5\u05d0\u05d8\u05d9\u05d1\u05d9\u05ea \u05dc-SSL(\u05ea\u05e2\u05e8 \u05d0\u05d5\u05e7\u05dd). \u05d1\u05de\u05e7\u05d5\u05dd \u05dc\u05e2\u05d1\u05d5\u05d3 \u05e2\u05dd \u05d4\u05d9\u05d9\u05e6\u05d5\u05d2\u05d9\u05dd \u05d4\u05de\u05d0\u05de\u05e8 \u05de\u05d0\u05de\u05df \u05de\u05d5\u05d3\u05dc \u05dc\u05d7\u05d6\u05d5\u05ea \u05d0\u05ea \u05d4\u05de\u05e1\u05e4\u05e8 \u05e9\u05dc \u05d4\u05d3\u05d5\u05d2\u05de\u05d0 \u05d1\u05d3\u05d0\u05d8\u05d4\u05e1\u05d8. \u05db\u05dc\u05d5\u05de\u05e8 \u05d0\u05dd \u05d9\u05e9 \u05dc\u05e0\u05d5 1000 \u05d3\u05d5\u05d2\u05de\u05d0\u05d5\u05ea \u05de\u05d4\u05d3\u05d0\u05d8\u05d4\u05e1\u05d8 \u05d9\u05e9 \u05dc\u05e0\u05d5 1000 \u05e7\u05d8\u05d2\u05d5\u05e8\u05d9\u05d5\u05ea \u05d5\u05de\u05d8\u05e8\u05ea\u05e0\u05d5 \u05dc\u05d7\u05d6\u05d5\u05ea \u05e7\u05d8\u05d2\u05d5\u05e8\u05d9\u05d4 \u05e9\u05dc \u05d3\u05d5\u05d2\u05de\u05d0 \u05de\u05d4\u05d9\u05d9\u05e6\u05d5\u05d2 \u05d4\u05dc\u05d8\u05e0\u05d8\u05d9 \u05e9\u05dc\u05d4 (\u05d4\u05de\u05d5\u05e4\u05e7 \u05e2\u05dc \u05d9\u05d3\u05d9 \u05d4\u05de\u05d5\u05d3\u05dc \u05d4\u05de\u05d0\u05d5\u05de\u05df). \u05db\u05dc\u05d5\u05de\u05e8 \u05d0\u05d7\u05e8\u05d9 \u05d4\u05e9\u05db\u05d1\u05d4 \u05d4\u05d0\u05d7\u05e8\u05d5\u05e0\u05d4 \u05e9\u05dc \u05d4\u05de\u05d5\u05d3\u05dc \u05de\u05d5\u05e1\u05d9\u05e4\u05d9\u05dd \u05e9\u05db\u05d1\u05d4 \u05e2\u05dd \u05de\u05d8\u05e8\u05d9\u05e6\u05d4 \u05d4\u05de\u05de\u05e4\u05d4 \u05d0\u05ea \u05d4\u05d9\u05d9\u05e6\u05d5\u05d2 \u05dc\u05e7\u05d8\u05d2\u05d5\u05e8\u05d9\u05d5\u05ea (\u05db\u05dc\u05d5\u05de\u05e8 \u05d4\u05de\u05e1\u05e4\u05e8\u05d9\u05dd \u05d4\u05e1\u05d9\u05d3\u05d5\u05e8\u05d9\u05d9\u05dd \u05e9\u05dc \u05d4\u05d3\u05d5\u05d2\u05de\u05d0\u05d5\u05ea). \u05d5\u05d1\u05e1\u05d5\u05e3 \u05e9\u05dc \u05dc\u05d5\u05e1 cross-entropy \u05d4\u05e1\u05d8\u05e0\u05d3\u05e8\u05d8\u05d9.\"\n\nSet-ParagraphText 6 \"\u05d0\u05d6 \u05d4\u05de\u05d0\u05de\u05e8 \u05de\u05d5\u05db\u05d9\u05d7 \u05e9\u05d4\u05e9\u05d9\u05d8\u05d4 \u05e2\u05d5\u05d1\u05d3\u05ea \u05dc\u05d0 \u05e8\u05e2 \u05dc\u05d3\u05d0\u05d8\u05d4\u05e1\u05d8\u05d9\u05dd \u05d9\u05d7\u05e1\u05d9\u05ea \u05dc\u05d0 \u05d2\u05d3\u05d5\u05dc\u05d9\u05dd (\u05de\u05e2\u05e0\u05d9\u05d9\u05df \u05d0\u05d9\u05da \u05d6\u05d4 \u05d9\u05e2\u05d1\u05d5\u05d3 \u05dc\u05d3\u05d0\u05d8\u05d4\u05e1\u05d8 \u05d1\u05d2\u05d5\u05d3\u05dc 10 \u05de\u05d9\u05dc\u05d9\u05d5\u05df). \u05db\u05de\u05d5\u05d1\u05df \u05d9\u05e9 \u05db\u05de\u05d4 \u05d8\u05e8\u05d9\u05e7\u05d9\u05dd \u05d1\u05d0\u05d9\u05de\u05d5\u05df \u05db\u05de\u05d5 soft labels \u05d0\u05d1\u05dc \u05d1\u05d2\u05d3\u05d5\u05dc \u05d4\u05e8\u05e2\u05d9\u05d5\u05df \u05d3\u05d9 \u05e0\u05d7\u05de\u05d3. \"\n\n# --- Paragraphs 7-9: entirely removed ---\n$startPara = $d.Paragraphs(7)\n$endPara = $d.Paragraphs(9)\n$range = $d.Range($startPara.Range.Start, $endPara.Range.End)\n$range.Delete()\n\n# --- Paragraph 10 (now last, index 7): replace URL text ---\nSet-ParagraphText 7 \"https://arxiv.org/pdf/2406.10743\"\n"}
